$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 546.7765096666667
$ws.Range("H2").Value = 1640.329529
$ws.Range("I2").Value = 0.6285526459909564
$ws.Range("J2").Value = 0.6285526459909564
$ws.Range("M2").Value = 9.423852333333334
$ws.Range("N2").Value = 28.271557
$ws.Range("O2").Value = 0.06654336290212845
$ws.Range("P2").Value = 0.06654336290212845
$ws.Range("Q2").Value = 5152.741086434074
$ws.Range("R2").Value = 46374.66977790666
$ws.Range("S2").Value = 0.04182600682526928
$ws.Range("T2").Value = 0.04182600682526928
$ws.Range("G3").Value = 546.7765096666667
$ws.Range("H3").Value = 1640.329529
$ws.Range("I3").Value = 0.6285526459909564
$ws.Range("J3").Value = 0.6285526459909564
$ws.Range("O3").Value = 0.3572423751649123
$ws.Range("P3").Value = 0.3572423751649123
$ws.Range("Q3").Value = 27662.82592352513
$ws.Range("R3").Value = 248965.4333117261
$ws.Range("S3").Value = 0.2245456401699996
$ws.Range("T3").Value = 0.2245456401699996
$ws.Range("G4").Value = 546.7765096666667
$ws.Range("H4").Value = 1640.329529
$ws.Range("I4").Value = 0.6285526459909564
$ws.Range("J4").Value = 0.6285526459909564
$ws.Range("M4").Value = 26.84076266666667
$ws.Range("N4").Value = 80.522288
$ws.Range("O4").Value = 0.1895270158659356
$ws.Range("P4").Value = 0.1895270158659356
$ws.Range("Q4").Value = 14675.89852767137
$ws.Range("R4").Value = 132083.0867490424
$ws.Range("S4").Value = 0.1191277073093038
$ws.Range("T4").Value = 0.1191277073093038
$ws.Range("G5").Value = 546.7765096666667
$ws.Range("H5").Value = 1640.329529
$ws.Range("I5").Value = 0.6285526459909564
$ws.Range("J5").Value = 0.6285526459909564
$ws.Range("M5").Value = 54.762539
$ws.Range("N5").Value = 164.287617
$ws.Range("O5").Value = 0.3866872460670236
$ws.Range("P5").Value = 0.3866872460670236
$ws.Range("Q5").Value = 29942.86993490471
$ws.Range("R5").Value = 269485.8294141424
$ws.Range("S5").Value = 0.2430532916863838
$ws.Range("T5").Value = 0.2430532916863838
$ws.Range("I6").Value = 0.1861770314550556
$ws.Range("J6").Value = 0.1861770314550556
$ws.Range("M6").Value = 9.423852333333334
$ws.Range("N6").Value = 28.271557
$ws.Range("O6").Value = 0.06654336290212845
$ws.Range("P6").Value = 0.06654336290212845
$ws.Range("Q6").Value = 1526.23975962484
$ws.Range("R6").Value = 13736.15783662355
$ws.Range("S6").Value = 0.01238884576815475
$ws.Range("T6").Value = 0.01238884576815475
$ws.Range("I7").Value = 0.1861770314550556
$ws.Range("J7").Value = 0.1861770314550556
$ws.Range("O7").Value = 0.3572423751649123
$ws.Range("P7").Value = 0.3572423751649123
$ws.Range("S7").Value = 0.06651032491815666
$ws.Range("T7").Value = 0.06651032491815666
$ws.Range("I8").Value = 0.1861770314550556
$ws.Range("J8").Value = 0.1861770314550556
$ws.Range("M8").Value = 26.84076266666667
$ws.Range("N8").Value = 80.522288
$ws.Range("O8").Value = 0.1895270158659356
$ws.Range("P8").Value = 0.1895270158659356
$ws.Range("Q8").Value = 4346.99501982017
$ws.Range("R8").Value = 39122.95517838153
$ws.Range("S8").Value = 0.03528557719445511
$ws.Range("T8").Value = 0.03528557719445512
$ws.Range("I9").Value = 0.1861770314550556
$ws.Range("J9").Value = 0.1861770314550556
$ws.Range("M9").Value = 54.762539
$ws.Range("N9").Value = 164.287617
$ws.Range("O9").Value = 0.3866872460670236
$ws.Range("P9").Value = 0.3866872460670236
$ws.Range("Q9").Value = 8869.065579918984
$ws.Range("R9").Value = 79821.59021927086
$ws.Range("S9").Value = 0.07199228357428909
$ws.Range("T9").Value = 0.07199228357428909
$ws.Range("G10").Value = 160.630483
$ws.Range("H10").Value = 481.891449
$ws.Range("I10").Value = 0.1846544489960017
$ws.Range("J10").Value = 0.1846544489960017
$ws.Range("M10").Value = 9.423852333333334
$ws.Range("N10").Value = 28.271557
$ws.Range("O10").Value = 0.06654336290212845
$ws.Range("P10").Value = 0.06654336290212845
$ws.Range("Q10").Value = 1513.757952024011
$ws.Range("R10").Value = 13623.82156821609
$ws.Range("S10").Value = 0.01228752801103351
$ws.Range("T10").Value = 0.01228752801103351
$ws.Range("G11").Value = 160.630483
$ws.Range("H11").Value = 481.891449
$ws.Range("I11").Value = 0.1846544489960017
$ws.Range("J11").Value = 0.1846544489960017
$ws.Range("O11").Value = 0.3572423751649123
$ws.Range("P11").Value = 0.3572423751649123
$ws.Range("Q11").Value = 8126.708098615401
$ws.Range("R11").Value = 73140.3728875386
$ws.Range("S11").Value = 0.0659663939440998
$ws.Range("T11").Value = 0.0659663939440998
$ws.Range("G12").Value = 160.630483
$ws.Range("H12").Value = 481.891449
$ws.Range("I12").Value = 0.1846544489960017
$ws.Range("J12").Value = 0.1846544489960017
$ws.Range("M12").Value = 26.84076266666667
$ws.Range("N12").Value = 80.522288
$ws.Range("O12").Value = 0.1895270158659356
$ws.Range("P12").Value = 0.1895270158659356
$ws.Range("Q12").Value = 4311.444671235035
$ws.Range("R12").Value = 38803.00204111531
$ws.Range("S12").Value = 0.03499700668458081
$ws.Range("T12").Value = 0.03499700668458081
$ws.Range("G13").Value = 160.630483
$ws.Range("H13").Value = 481.891449
$ws.Range("I13").Value = 0.1846544489960017
$ws.Range("J13").Value = 0.1846544489960017
$ws.Range("M13").Value = 54.762539
$ws.Range("N13").Value = 164.287617
$ws.Range("O13").Value = 0.3866872460670236
$ws.Range("P13").Value = 0.3866872460670236
$ws.Range("Q13").Value = 8796.533089876337
$ws.Range("R13").Value = 79168.79780888701
$ws.Range("S13").Value = 0.07140352035628758
$ws.Range("T13").Value = 0.07140352035628758
$ws.Range("G14").Value = 0.5357470000000001
$ws.Range("H14").Value = 1.607241
$ws.Range("I14").Value = 0.0006158735579862568
$ws.Range("J14").Value = 0.0006158735579862568
$ws.Range("M14").Value = 9.423852333333334
$ws.Range("N14").Value = 28.271557
$ws.Range("O14").Value = 0.06654336290212845
$ws.Range("P14").Value = 0.06654336290212845
$ws.Range("Q14").Value = 5.048800616026335
$ws.Range("R14").Value = 45.43920554423701
$ws.Range("S14").Value = 0.00004098229767090454
$ws.Range("T14").Value = 0.00004098229767090454
$ws.Range("G15").Value = 0.5357470000000001
$ws.Range("H15").Value = 1.607241
$ws.Range("I15").Value = 0.0006158735579862568
$ws.Range("J15").Value = 0.0006158735579862568
$ws.Range("O15").Value = 0.3572423751649123
$ws.Range("P15").Value = 0.3572423751649123
$ws.Range("Q15").Value = 27.10481474247267
$ws.Range("R15").Value = 243.943332682254
$ws.Range("S15").Value = 0.0002200161326562757
$ws.Range("T15").Value = 0.0002200161326562757
$ws.Range("G16").Value = 0.5357470000000001
$ws.Range("H16").Value = 1.607241
$ws.Range("I16").Value = 0.0006158735579862568
$ws.Range("J16").Value = 0.0006158735579862568
$ws.Range("M16").Value = 26.84076266666667
$ws.Range("N16").Value = 80.522288
$ws.Range("O16").Value = 0.1895270158659356
$ws.Range("P16").Value = 0.1895270158659356
$ws.Range("Q16").Value = 14.37985807637867
$ws.Range("R16").Value = 129.418722687408
$ws.Range("S16").Value = 0.0001167246775958715
$ws.Range("T16").Value = 0.0001167246775958715
$ws.Range("G17").Value = 0.5357470000000001
$ws.Range("H17").Value = 1.607241
$ws.Range("I17").Value = 0.0006158735579862568
$ws.Range("J17").Value = 0.0006158735579862568
$ws.Range("M17").Value = 54.762539
$ws.Range("N17").Value = 164.287617
$ws.Range("O17").Value = 0.3866872460670236
$ws.Range("P17").Value = 0.3866872460670236
$ws.Range("Q17").Value = 29.338865981633
$ws.Range("R17").Value = 264.049793834697
$ws.Range("S17").Value = 0.000238150450063205
$ws.Range("T17").Value = 0.000238150450063205

Write-Host "Applied all updates"
